# Auto-generated: apply scheduled-runner market price refresh to Belias_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1772.2727
$ws.Range("I40").Value = 1551.55
$ws.Range("J40").Value = 2111.8462
$ws.Range("K40").Value = 1551.55
$ws.Range("L40").Value = 2111.8462
$ws.Range("M40").Value = -1376.55
$ws.Range("N40").Value = -2461.8462

$ws.Range("H101").Value = 450.375
$ws.Range("I101").Value = 376
$ws.Range("K101").Value = 1128
$ws.Range("M101").Value = 494

$ws.Range("H103").Value = 3720.647
$ws.Range("J103").Value = 749.8
$ws.Range("L103").Value = 2249.4
$ws.Range("N103").Value = -3421.4

$ws.Range("H125").Value = 2372.6667
$ws.Range("I125").Value = 5000
$ws.Range("J125").Value = 1847.2
$ws.Range("K125").Value = 45000
$ws.Range("L125").Value = 16624.8
$ws.Range("M125").Value = -42540
$ws.Range("N125").Value = -21544.8

$ws.Range("H132").Value = 4887.9062
$ws.Range("I132").Value = 1738.862
$ws.Range("J132").Value = 35328.668
$ws.Range("K132").Value = 5216.586
$ws.Range("L132").Value = 105986.004
$ws.Range("M132").Value = -2686.586
$ws.Range("N132").Value = -111046.004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2410.182
$ws.Range("I61").Value = 2279.111
$ws.Range("K61").Value = 2279.111
$ws.Range("M61").Value = -2067.111

$ws.Range("H122").Value = 1876.4706
$ws.Range("J122").Value = 1788.4615
$ws.Range("L122").Value = 5365.3845
$ws.Range("N122").Value = -10265.3845

$ws.Range("H136").Value = 2410.182
$ws.Range("I136").Value = 2279.111
$ws.Range("K136").Value = 6837.333
$ws.Range("M136").Value = -4287.333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H43").Value = 89342
$ws.Range("J43").Value = 89342
$ws.Range("L43").Value = 89342
$ws.Range("N43").Value = -89704

$ws.Range("H92").Value = 29351
$ws.Range("J92").Value = 29351
$ws.Range("L92").Value = 29351
$ws.Range("N92").Value = -34343

$ws.Range("H134").Value = 2112033.5
$ws.Range("I134").Value = 3083818.8
$ws.Range("K134").Value = 9251456.399999999
$ws.Range("M134").Value = -9248921.399999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 2591367
$ws.Range("I86").Value = 5854788
$ws.Range("J86").Value = 7825.25
$ws.Range("K86").Value = 5854788
$ws.Range("L86").Value = 7825.25
$ws.Range("M86").Value = -5853665
$ws.Range("N86").Value = -10071.25

$ws.Range("H89").Value = 2591367
$ws.Range("I89").Value = 5854788
$ws.Range("J89").Value = 7825.25
$ws.Range("K89").Value = 29273940
$ws.Range("L89").Value = 39126.25
$ws.Range("M89").Value = -29268324
$ws.Range("N89").Value = -50358.25

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 50546.9
$ws.Range("I2").Value = 30
$ws.Range("J2").Value = 84224.836
$ws.Range("K2").Value = 180
$ws.Range("L2").Value = 505349.0159999999
$ws.Range("M2").Value = -67
$ws.Range("N2").Value = -505575.0159999999

$ws.Range("H38").Value = 192.58333
$ws.Range("I38").Value = 60.25
$ws.Range("K38").Value = 180.75
$ws.Range("M38").Value = 166.25

$ws.Range("H113").Value = 1317999.5
$ws.Range("I113").Value = 2525741.5
$ws.Range("J113").Value = 462.81818
$ws.Range("K113").Value = 7577224.5
$ws.Range("L113").Value = 1388.45454
$ws.Range("M113").Value = -7575054.5
$ws.Range("N113").Value = -5728.45454

$ws.Range("H131").Value = 888.86
$ws.Range("I131").Value = 585.6667
$ws.Range("J131").Value = 908.21277
$ws.Range("K131").Value = 1757.0001
$ws.Range("L131").Value = 2724.63831
$ws.Range("M131").Value = 3282.9999
$ws.Range("N131").Value = -12804.63831

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4630631
$ws.Range("I102").Value = 11905480
$ws.Range("J102").Value = 1181.2727
$ws.Range("K102").Value = 11905480
$ws.Range("L102").Value = 1181.2727
$ws.Range("M102").Value = -11903858
$ws.Range("N102").Value = -4425.2727

$ws.Range("H122").Value = 101087.4
$ws.Range("I122").Value = 126209.25
$ws.Range("K122").Value = 378627.75
$ws.Range("M122").Value = -376177.75

$ws.Range("H132").Value = 2566879
$ws.Range("I132").Value = 2978.5833
$ws.Range("K132").Value = 8935.749899999999
$ws.Range("M132").Value = -6405.749899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5276.6665
$ws.Range("J7").Value = 4938
$ws.Range("L7").Value = 4938
$ws.Range("N7").Value = -5162

$ws.Range("H16").Value = 1665
$ws.Range("I16").Value = 1727.2858
$ws.Range("J16").Value = 1447
$ws.Range("K16").Value = 1727.2858
$ws.Range("L16").Value = 1447
$ws.Range("M16").Value = -1557.2858
$ws.Range("N16").Value = -1787

$ws.Range("H93").Value = 65754.91
$ws.Range("I93").Value = 1483.3334
$ws.Range("K93").Value = 1483.3334
$ws.Range("M93").Value = -235.3334

$ws.Range("H126").Value = 5276.6665
$ws.Range("J126").Value = 4938
$ws.Range("L126").Value = 14814
$ws.Range("N126").Value = -19754

$ws.Range("H132").Value = 3439.0344
$ws.Range("I132").Value = 3230.95
$ws.Range("J132").Value = 3901.4443
$ws.Range("K132").Value = 9692.849999999999
$ws.Range("L132").Value = 11704.3329
$ws.Range("M132").Value = -7162.849999999999
$ws.Range("N132").Value = -16764.3329

$ws.Range("H136").Value = 1786.3043
$ws.Range("I136").Value = 1240
$ws.Range("K136").Value = 3720
$ws.Range("M136").Value = -1170

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 56449.684
$ws.Range("I122").Value = 1653.091
$ws.Range("K122").Value = 4959.272999999999
$ws.Range("M122").Value = -2509.272999999999

$ws.Range("H126").Value = 1800.3636
$ws.Range("I126").Value = 1640.8
$ws.Range("J126").Value = 1933.3334
$ws.Range("K126").Value = 4922.4
$ws.Range("L126").Value = 5800.0002
$ws.Range("M126").Value = -2452.4
$ws.Range("N126").Value = -10740.0002

$ws.Range("H132").Value = 2835.9268
$ws.Range("I132").Value = 2502.4482
$ws.Range("J132").Value = 3641.8333
$ws.Range("K132").Value = 7507.344599999999
$ws.Range("L132").Value = 10925.4999
$ws.Range("M132").Value = -4977.344599999999
$ws.Range("N132").Value = -15985.4999
